$wb = $excel.ActiveWorkbook

# Map of sheet name -> hashtable of cell address -> new value
$updates = @{
    "展览" = @{
        "F2"  = 709
        "F3"  = 13573
        "F4"  = 13816
        "F8"  = 5752
        "F16" = 2125
        "F18" = 1745
        "F28" = 59
        "F31" = 1745
        "F33" = 1292
        "F34" = 87
        "F36" = 4550
        "F37" = 4659
        "F42" = 3235
        "F46" = 75
        "F49" = 257
    }
    "本地生活" = @{
        "F2" = 7322
        "F3" = 197
        "F4" = 613
    }
    "全部类型" = @{
        "F2"  = 709
        "F3"  = 197
        "F4"  = 613
        "F6"  = 13573
        "F7"  = 13816
        "F10" = 5752
        "F18" = 2125
        "F20" = 1745
        "F25" = 59
        "F28" = 1745
        "F31" = 1292
        "F33" = 87
        "F34" = 4550
        "F35" = 4659
        "F41" = 3235
        "F45" = 75
        "F48" = 257
    }
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $cellMap = $updates[$sheetName]
    foreach ($addr in $cellMap.Keys) {
        $ws.Range($addr).Value = $cellMap[$addr]
    }
}
